$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.615.71"
$ws.Range("E2").Value = "  +4.64%  "
$ws.Range("D3").Value = "3.089.32"
$ws.Range("E3").Value = "  +2.66%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.35"
$ws.Range("E5").Value = "  +2.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.74"
$ws.Range("E6").Value = "  +8.41%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.084.28"
$ws.Range("E8").Value = "  +2.72%  "
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.15"
$ws.Range("E10").Value = "  +17.14%  "
$ws.Range("E11").Value = "  +3.90%  "
$ws.Range("E12").Value = "  +3.84%  "
$ws.Range("E13").Value = "  +3.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.25"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "3.597.09"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "64.636.62"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "3.092.83"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "481.62"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.78"
$ws.Range("E21").Value = "  +4.05%  "
$ws.Range("E22").Value = "  +0.96%  "
$ws.Range("E23").Value = "  +8.53%  "
$ws.Range("E24").Value = "  +11.65%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.82"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  +3.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  +5.70%  "
$ws.Range("E29").Value = "  +7.64%  "
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.00"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  +3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.45"
$ws.Range("E33").Value = "  +5.45%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.21"
$ws.Range("E35").Value = "  +6.02%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "465.32"
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0406"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("E39").Value = "  +3.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.85"
$ws.Range("E40").Value = "  +16.86%  "
$ws.Range("D41").Value = "3.007.64"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.25"
$ws.Range("E42").Value = "  +2.13%  "
$ws.Range("E43").Value = "  -1.39%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "28.19"
$ws.Range("E44").Value = "  +6.75%  "
$ws.Range("E45").Value = "  +6.78%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +7.84%  "
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.94"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("D50").Value = "0.0₃0515"
$ws.Range("E50").Value = "  +5.57%  "
$ws.Range("E51").Value = "  +2.46%  "
